# Updated symbol list on Fri Dec 23 19:59:07 UTC 2022 with GitHub Actions
#
# This script applies updated crypto price / ordering data to Sheet1.
# All data cells in the sheet are stored as text (inline strings) even
# when they look numeric (e.g. "246.24"), so for any new value that looks
# like a number we force the cell to Text format before assigning it and
# then reset the style back to Normal so no stray formatting is left
# behind (matching the original workbook, where these cells carry no
# explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, [string]$Address, [string]$Text)
    $range = $Worksheet.Range($Address)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = "Normal"
}

function Set-PlainValue {
    param($Worksheet, [string]$Address, [string]$Text)
    $Worksheet.Range($Address).Value = $Text
}

# --- Price (column D) updates -------------------------------------------------
Set-TextValue $ws "D2"  "246.13"
Set-TextValue $ws "D3"  "22.23"
Set-TextValue $ws "D4"  "5.352"
Set-TextValue $ws "D5"  "0.05863"
Set-TextValue $ws "D7"  "6.377"
Set-TextValue $ws "D8"  "0.8127"
Set-TextValue $ws "D9"  "1.000"
Set-TextValue $ws "D10" "0.1421"
Set-TextValue $ws "D11" "0.03534"
Set-TextValue $ws "D12" "0.07365"
Set-TextValue $ws "D13" "0.03004"
Set-TextValue $ws "D14" "4.192"
Set-TextValue $ws "D15" "0.09396"
Set-TextValue $ws "D16" "0.001591"
Set-TextValue $ws "D17" "0.04836"

# Row 18 (One/ONE) - "Worstin24h" tag added to the volume label
Set-PlainValue $ws "E18" "17OneONEWorstin24h"

Set-TextValue $ws "D19" "0.006218"
Set-TextValue $ws "D20" "0.004081"
Set-TextValue $ws "D21" "0.0009848"
Set-TextValue $ws "D22" "0.00010000"
Set-TextValue $ws "D23" "3.689"
Set-TextValue $ws "D27" "0.0002472"

# --- Rows 41-43 rotate: KickToken -> row41, BKEXToken -> row42, CEJI -> row43 --
Set-PlainValue $ws "B41" "KickToken"
Set-PlainValue $ws "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue   $ws "D41" "0.006470"
Set-PlainValue $ws "E41" "40KickTokenKICK"

Set-PlainValue $ws "B42" "BKEXToken"
Set-PlainValue $ws "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue   $ws "D42" "0.1075"
Set-PlainValue $ws "E42" "41BKEXTokenBKK"

Set-PlainValue $ws "B43" "CEJI"
Set-PlainValue $ws "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue   $ws "D43" "0.003000"
Set-PlainValue $ws "E43" "42CEJICEJI"

# --- Remaining price (column D) updates ---------------------------------------
Set-TextValue $ws "D44" "0.005220"
Set-TextValue $ws "D45" "0.00005651"
Set-TextValue $ws "D47" "0.7221"
Set-TextValue $ws "D48" "0.08106"
